# Brainstorm.docx edit:
# After the last paragraph ("Je kan doormiddel van de mogelijke bruggen een
# amino verplaatsen van plek om een brug te maken.") add a series of new
# paragraphs (two blank lines, "Hillclimb bekijken", a blank line, then
# "Liesbeth:", "Ander algoritme", "Wat gaat ze doen?") right before the
# section break that closes the document body.

$d = $word.ActiveDocument

# Move to the very end of the document's main story.
$r = $d.Content
$r.Collapse(1)  # wdCollapseEnd

function New-Paragraph {
    # Appends one more (initially empty) paragraph right after the current
    # end of the document, mirroring pressing Enter at the end of the text.
    $rng = $d.Content
    $rng.Collapse(0)  # wdCollapseEnd
    $rng.InsertParagraphAfter()
}

function New-TextParagraph([string]$text) {
    New-Paragraph
    $rng = $d.Content
    $rng.Collapse(0)  # wdCollapseEnd
    $rng.InsertAfter($text)
}

# Two blank paragraphs.
New-Paragraph
New-Paragraph

# "Hillclimb bekijken"
New-TextParagraph "Hillclimb bekijken"

# Blank paragraph.
New-Paragraph

# Liesbeth notes.
New-TextParagraph "Liesbeth:"
New-TextParagraph "Ander algoritme"
New-TextParagraph "Wat gaat ze doen?"
